$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'folgen',
    'saufen',
    'schwächen',
    'wüten',
    'filmen',
    'schmecken',
    'mauern',
    'biegen',
    'scheinen',
    'drehen',
    'fahren',
    'graben',
    'sorgen',
    'suchen',
    'fischen',
    'tropfen',
    'klettern',
    'heilen',
    'ärgern',
    'malen',
    'wenden',
    'weichen',
    'runden',
    'schulden',
    'pflanzen',
    'streichen',
    'seufzen',
    'lesen',
    'kosten',
    'bauen',
    'flüchten',
    'klingen',
    'stecken',
    'wirken',
    'achten',
    'greifen',
    'schwingen',
    'rufen',
    'warnen',
    'äußern',
    'töten',
    'liegen',
    'geben',
    'trennen',
    'kichern',
    'arten',
    'altern',
    'decken',
    'formen',
    'planen',
    'treiben',
    'loben',
    'trauen',
    'kümmern',
    'dienen',
    'sperren',
    'spüren',
    'gelten',
    'sichern',
    'heben',
    'schrecken',
    'feiern',
    'lockern',
    'helfen',
    'hören',
    'lügen',
    'sprengen',
    'bergen',
    'rasen',
    'schwören',
    'sterben',
    'bluten',
    'schlucken',
    'münzen',
    'mögen',
    'irren',
    'knarren',
    'jubeln',
    'pfeifen',
    'ändern',
    'fällen',
    'grüßen',
    'stammen',
    'schenken',
    'schreiten',
    'fangen',
    'brauchen',
    'quälen',
    'jagen',
    'zünden',
    'fallen',
    'scheitern',
    'wehtun',
    'werden',
    'betteln',
    'siegen',
    'spielen',
    'bitten',
    'kehren',
    'stehlen',
    'räumen',
    'heulen',
    'reizen',
    'ehren',
    'hauen',
    'machen',
    'wundern',
    'zögern',
    'führen',
    'liefern',
    'platzen',
    'werfen',
    'spinnen',
    'gründen',
    'enden',
    'sinken',
    'erben',
    'fließen',
    'boxen',
    'backen',
    'wachsen',
    'zielen',
    'tollen',
    'bellen',
    'flehen',
    'freuen',
    'zeigen',
    'dringen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}

